# Updated Global_M2 for easier usage.
# Append three new monthly FX rows (212-214) and correct the low/close
# values of the last existing row (211) on the Iraq FX sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the existing last row (211): low/close were placeholders ---
$ws.Cells.Item(211, 5).Value = 1308   # E211 low:  1311 -> 1308
$ws.Cells.Item(211, 6).Value = 1308   # F211 close: 1459 -> 1308

# --- New row 212 ---
$ws.Cells.Item(211, 1).Copy()
$ws.Cells.Item(212, 1).PasteSpecial(-4122)   # xlPasteFormats, carries date style (s=2)
$ws.Cells.Item(212, 1).Value = 45047.33333333334
$ws.Cells.Item(212, 2).Value = "FX_IDC:USDIQD"
$ws.Cells.Item(212, 3).Value = 1309
$ws.Cells.Item(212, 4).Value = 1309
$ws.Cells.Item(212, 5).Value = 1307
$ws.Cells.Item(212, 6).Value = 1309
$ws.Cells.Item(212, 7).Value = 0

# --- New row 213 ---
$ws.Cells.Item(211, 1).Copy()
$ws.Cells.Item(213, 1).PasteSpecial(-4122)
$ws.Cells.Item(213, 1).Value = 45078.33333333334
$ws.Cells.Item(213, 2).Value = "FX_IDC:USDIQD"
$ws.Cells.Item(213, 3).Value = 1309
$ws.Cells.Item(213, 4).Value = 1309
$ws.Cells.Item(213, 5).Value = 1307
$ws.Cells.Item(213, 6).Value = 1307
$ws.Cells.Item(213, 7).Value = 0

# --- New row 214 ---
$ws.Cells.Item(211, 1).Copy()
$ws.Cells.Item(214, 1).PasteSpecial(-4122)
$ws.Cells.Item(214, 1).Value = 45110.33333333334
$ws.Cells.Item(214, 2).Value = "FX_IDC:USDIQD"
$ws.Cells.Item(214, 3).Value = 1308
$ws.Cells.Item(214, 4).Value = 1308
$ws.Cells.Item(214, 5).Value = 1307
$ws.Cells.Item(214, 6).Value = 1307
$ws.Cells.Item(214, 7).Value = 0

$excel.CutCopyMode = 0
